$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19, column E: the error type changed from ValueError to TypeError
# (same message text, different exception type)
$ws.Range("E19").Value2 = 'TypeError("lista invalida, não suportado (String)")'

# Copy formatting of the last existing data row (20) down into the new row (21)
$ws.Range("A20:E20").Copy()
$ws.Range("A21:E21").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Fill in the new row of test data
$ws.Range("A21").Value2 = "escola/aluno.py"
$ws.Range("B21").Value2 = "calcular_media"
$ws.Range("D21").Value2 = '"ola"'
$ws.Range("C21").Value2 = "enviando uma string"
$ws.Range("E21").Value2 = 'TypeError("não é uma lista")'

# Update the selected cell to match the saved view state
$ws.Range("E24").Select()
